# "Generate Report for Handback"
#
# Before: every file is still "Ready for handoff" and the per-language
# sheets (zh-cn / de-de) have no Latest Target File / Latest Handback File
# recorded yet, and the Latest Handback DateTime is the zero-date sentinel.
#
# After: the two real source files (ffb66c72-...md and fffff3187753-...md)
# have been handed back and are in sync with en-US, so:
#   - Status everywhere flips to "Handed back: in sync with en-US"
#   - Latest Target File / Latest Handback File get filled in (with the
#     same hyperlinks as the existing Source File Name / Latest Handoff
#     File columns)
#   - Latest Handback DateTime records when the handback happened
#
# The ".localization-config" row (row 4) is excluded from localization and
# is left untouched.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

$mdFile  = "ffb66c72-7a23-47c6-82db-d0759d329a7f.md"
$mdUrl   = "https://github.com/OpenLocalizationTest/oltest/blob/451ce9e589f699dbd6bb724e10671ebd88ff894c/e2e/ffb66c72-7a23-47c6-82db-d0759d329a7f.md"

$zhXlf   = "ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.zh-cn.xlf"
$zhUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8cc5869b13bb3c1cf0daa94d4c7cb94e2f686339/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/yuwzho/ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.zh-cn.xlf"

$deXlf   = "ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.de-de.xlf"
$deUrl   = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d988695e908499e2428d8a3b5dd47463b9eba9d6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/yuwzho/ffb66c72-7a23-47c6-82db-d0759d329a7f.41d9e9e7e2c2c44afa6210287a70fd273598c3bf.de-de.xlf"

$zhHandback = "2016-01-25 08:36:31"
$deHandback = "2016-01-25 08:36:52"

# ---------------------------------------------------------------------
# Overview sheet: flip the two per-language status columns for the two
# real files (row 2 = ffb66c72..md, row 3 = fffff3187753..md). Row 4 is
# the .localization-config row ("Not to be localized") and stays as-is.
# ---------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")
$ov.Range("B2").Value = $newStatus
$ov.Range("C2").Value = $newStatus
$ov.Range("B3").Value = $newStatus
$ov.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Helper: populate one language sheet (zh-cn / de-de) with the handback
# information for its two real-file rows (2 and 3).
# ---------------------------------------------------------------------
function Update-LanguageSheet($sheetName, $xlfFile, $xlfUrl, $handbackTime) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in 2, 3) {
        # Status -> handed back
        $ws.Range("B$row").Value = $newStatus

        # Latest Target File == the source md file that was localized,
        # same file/link as column A (Source File Name).
        $ws.Range("E$row").Value = $mdFile
        $ws.Hyperlinks.Add($ws.Range("E$row"), $mdUrl, "", "", $mdFile) | Out-Null
        # Match the workbook's existing custom HyperLink look (underline,
        # cornflower blue) instead of Excel's auto-generated theme style.
        $ws.Range("E$row").Font.Underline = 2
        $ws.Range("E$row").Font.Color = 15570276

        # Latest Handback File == the handed-off xlf file came back,
        # same file/link as column C (Latest Handoff File).
        $ws.Range("F$row").Value = $xlfFile
        $ws.Hyperlinks.Add($ws.Range("F$row"), $xlfUrl, "", "", $xlfFile) | Out-Null
        $ws.Range("F$row").Font.Underline = 2
        $ws.Range("F$row").Font.Color = 15570276

        # Latest Handback DateTime -- when the handback happened.
        $ws.Range("G$row").Value = $handbackTime
    }
}

Update-LanguageSheet "zh-cn" $zhXlf $zhUrl $zhHandback
Update-LanguageSheet "de-de" $deXlf $deUrl $deHandback
